# The Pharmacy deposit test case (TC017 / Pharmacy\TC006PharmacyDeposit&Return.py)
# row has been deprecated - remove its entire row from the test-case list.
# Everything below shifts up by one row automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the whole row before deleting it (mirrors the author selecting the
# row in the UI prior to removing it) so the saved selection lands on the
# row that now occupies position 17 (Pharmacy\TC007CreatePharmacyGoodsReceipt.py).
$ws.Rows(17).EntireRow.Select()
$ws.Rows(17).EntireRow.Delete()

$ws.Range("A17").Select()
